$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: value 0, bold/centered/bordered style
$rng = $ws.Range("B1")
$rng.Value = 0
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.VerticalAlignment = -4160     # xlTop
$rng.Borders.LineStyle = 1         # xlContinuous
$rng.Borders.Weight = 2            # xlThin

# Reuse the exact same style for A2 via a format-only copy/paste so the
# engine reuses the cellXf record instead of generating a new one.
$rng.Copy()
$ws.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0

# B2: plain shared-string text value, no special style
$ws.Range("B2").Value = "disconnected_elements"
